# Commit: "Added more struggle data and repaired falling data"
# The sensor data rows (C2:H21) are shifted down by one row (row N's
# existing values move to row N+1), and a brand-new row of data is
# written into row 2 (the previous last row, 21, is dropped since the
# sheet's dimensions/timestamps/labels stay the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2..20 down into rows 3..21.
# Walk bottom-up so we don't clobber a value before it's copied.
for ($r = 20; $r -ge 2; $r--) {
    for ($c = 3; $c -le 8; $c++) {
        $srcVal = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 1, $c).Value = $srcVal
    }
}

# Write the new values for row 2.
$ws.Cells.Item(2, 3).Value = 1.442139625549314
$ws.Cells.Item(2, 4).Value = 0.4175686836242678
$ws.Cells.Item(2, 5).Value = -0.7117971777915959
$ws.Cells.Item(2, 6).Value = -0.1186605766415596
$ws.Cells.Item(2, 7).Value = -0.207236036658287
$ws.Cells.Item(2, 8).Value = -0.081550508737564
